# Ticket creation page added. Login information extraction remaining
#
# This script reproduces, via the Excel COM object model, the edits that were
# made to the weekly plan: the calendar header row is extended by one more
# day (Wed, column AE), that day's date is filled in, and the four task
# "swim-lanes" (Admin / Project Manager / Developer / Submitter) have their
# milestone labels shifted right to make room for a new "(Assign Tickets)" /
# "(Create Ticket)" pair of columns used for the ticket-creation page.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 4 — weekday header strip. Extend the repeating Tue..Mon pattern by one
# more column (AE4 = "Wed"), continuing straight on from AD4 = "Tue" (which
# used to hold the special "Tuesday" label).
# ---------------------------------------------------------------------------
$ws.Range("AD4").Value = "Tue"
$ws.Range("AE4").Value = "Wed"

# ---------------------------------------------------------------------------
# Row 5 — date strip. Fill in the date for the newly added column.
# ---------------------------------------------------------------------------
$ws.Range("AE5").Value = 44188

# ---------------------------------------------------------------------------
# Row 7 (Admin lane) — shift "Landing Page" / "User Management Tab" /
# "Project Management + Tickets" to the right to make room.
# ---------------------------------------------------------------------------
$ws.Range("E7").Clear()
$ws.Range("N7").Clear()

$ws.Range("G7").Value = "Landing Page"
$ws.Range("G7").Font.Bold = $true

$ws.Range("I7").Value = "User Management Tab"

$ws.Range("P7").Value = "Project Management + Tickets"
$ws.Range("P7").Font.Bold = $true

$ws.Range("Y7").Value = "Test + UI"

# ---------------------------------------------------------------------------
# Row 9 (Project Manager lane) — "(Assign Tickets)" moves from D9 to F9 (and
# loses its bold styling), "Front + Back End (Test + UI)" / "Deploy" shift
# from AC9/AD9 to AD9/AE9.
# ---------------------------------------------------------------------------
$ws.Range("D9").Clear()
$ws.Range("AC9").Clear()

$ws.Range("F9").Value = "(Assign Tickets)"

$ws.Range("AD9").Value = "Front + Back End (Test + UI)"
$ws.Range("AD9").Font.Bold = $true

$ws.Range("AE9").Value = "Deploy"
$ws.Range("AE9").Font.Bold = $true

# ---------------------------------------------------------------------------
# Row 11 (Developer lane) — "(Create Ticket)" moves from B11 to D11 (and
# loses its bold styling), "* Front end + Back End" moves from D11 to F11
# (gains bold styling), "Test + UI" stays in Y11.
# ---------------------------------------------------------------------------
$ws.Range("B11").Clear()

$ws.Range("D11").Value = "(Create Ticket)"
$ws.Range("D11").Font.Bold = $false

$ws.Range("F11").Value = "* Front end + Back End"
$ws.Range("F11").Font.Bold = $true

$ws.Range("Y11").Value = "Test + UI"

# ---------------------------------------------------------------------------
# Row 13 (Submitter lane) — "Front end + Back End" moves from B13 to D13.
# ---------------------------------------------------------------------------
$ws.Range("B13").Clear()

$ws.Range("D13").Value = "Front end + Back End"
$ws.Range("D13").Font.Bold = $true

$ws.Range("Y13").Value = "Test + UI"

# ---------------------------------------------------------------------------
# Column width tweaks that came along with the relayout (Excel's own
# autofit recompute for the columns whose content changed width).
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 8.0533854166667
$ws.Columns("H").ColumnWidth = 12.2760416666667
$ws.Columns("O").ColumnWidth = 13.9440104166667

# ---------------------------------------------------------------------------
# Selection cursor moved from H11 to G12.
# ---------------------------------------------------------------------------
$ws.Range("G12").Select()
